$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two new rows at 826-827, shifting existing rows 826-942 down to 828-944
$ws.Range("A826:R827").Insert()

# Populate new row 826 (Escarola / Primera, new sampling date 44776)
$ws.Cells.Item(826,1).Value = 1
$ws.Cells.Item(826,2).Value = "Agrícola del Norte S.A. de Arica"
$ws.Cells.Item(826,3).Value = "Arica y Parinacota"
$ws.Cells.Item(826,4).Value = 44776
$ws.Cells.Item(826,5).Value = 15
$ws.Cells.Item(826,6).Value = 100112033
$ws.Cells.Item(826,7).Value = "Lechuga"
$ws.Cells.Item(826,8).Value = "Escarola"
$ws.Cells.Item(826,9).Value = "Primera"
$ws.Cells.Item(826,10).Value = 120
$ws.Cells.Item(826,11).Value = 5000
$ws.Cells.Item(826,12).Value = 6000
$ws.Cells.Item(826,13).Value = 5500
$ws.Cells.Item(826,14).Value = "`$/caja 12 unidades"
$ws.Cells.Item(826,15).Value = "Región de Arica y Parinacota"
$ws.Cells.Item(826,16).Value = 458
$ws.Cells.Item(826,17).Value = 12
$ws.Cells.Item(826,18).Value = "Hortaliza"

# Populate new row 827 (Escarola / Segunda, new sampling date 44776)
$ws.Cells.Item(827,1).Value = 1
$ws.Cells.Item(827,2).Value = "Agrícola del Norte S.A. de Arica"
$ws.Cells.Item(827,3).Value = "Arica y Parinacota"
$ws.Cells.Item(827,4).Value = 44776
$ws.Cells.Item(827,5).Value = 15
$ws.Cells.Item(827,6).Value = 100112033
$ws.Cells.Item(827,7).Value = "Lechuga"
$ws.Cells.Item(827,8).Value = "Escarola"
$ws.Cells.Item(827,9).Value = "Segunda"
$ws.Cells.Item(827,10).Value = 120
$ws.Cells.Item(827,11).Value = 5000
$ws.Cells.Item(827,12).Value = 6000
$ws.Cells.Item(827,13).Value = 5500
$ws.Cells.Item(827,14).Value = "`$/caja 18 unidades"
$ws.Cells.Item(827,15).Value = "Región de Arica y Parinacota"
$ws.Cells.Item(827,16).Value = 306
$ws.Cells.Item(827,17).Value = 18
$ws.Cells.Item(827,18).Value = "Hortaliza"
